# Trade #73 closed at 2026-02-17 12:57:15 - unknown UNKNOWN +0.000%
# Applies the updated Summary / Strategy Status metrics and appends the
# new closed MarketMaking trade row (#73) to both the "All Trades" and
# "MarketMaking" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.26           # Current Capital
$summary.Range("B4").Value = 0.25              # Total P&L $
$summary.Range("B5").Value = 0.07000000000000001  # Total P&L %
$summary.Range("B6").Value = 73                # Total Trades
$summary.Range("B7").Value = 33                # Winning Trades
$summary.Range("B9").Value = 45.21             # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet (MarketMaking row, row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.26             # Capital
$status.Range("D4").Value = 73                 # Trades
$status.Range("E4").Value = 0.25               # P&L $
$status.Range("F4").Value = 0.26               # P&L %
$status.Range("G4").Value = 45.21              # Win Rate %

# ---------------------------------------------------------------------
# New trade row data shared by "All Trades" and "MarketMaking" sheets
# ---------------------------------------------------------------------
$newRow = @(73, "2026-02-17", "12:57:08", "MarketMaking", "UP", 0.57, 0.62, "CLOSED", 8.7719, 0.05, 100.26, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.11)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 74
    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $col = $i + 1
        $cell = $ws.Cells.Item($row, $col)
        # Column B ("2026-02-17") looks like a date: force it to stay the
        # plain text the source data uses (matching every other row),
        # instead of being auto-converted to a date serial number.
        if ($col -eq 2) {
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$i]
            $cell.ClearFormats()
        } else {
            $cell.Value = $newRow[$i]
        }
    }
}

$wb.Save()
